$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 3743.1428
$ws.Range("J32").Value = 4117
$ws.Range("L32").Value = 4117
$ws.Range("N32").Value = -4769
# Row 33
$ws.Range("H33").Value = 166.07143
$ws.Range("I33").Value = 166.07143
$ws.Range("K33").Value = 166.07143
$ws.Range("M33").Value = 62.92857000000001
# Row 43
$ws.Range("H43").Value = 2493
$ws.Range("J43").Value = 2493
$ws.Range("L43").Value = 2493
$ws.Range("N43").Value = -2631
# Row 64
$ws.Range("H64").Value = 4673.875
$ws.Range("I64").Value = 3798.5
$ws.Range("K64").Value = 3798.5
$ws.Range("M64").Value = -3550.5
# Row 67
$ws.Range("H67").Value = 4673.875
$ws.Range("I67").Value = 3798.5
$ws.Range("K67").Value = 3798.5
$ws.Range("M67").Value = -2940.5
# Row 98
$ws.Range("H98").Value = 908.875
$ws.Range("I98").Value = 854.8
$ws.Range("K98").Value = 854.8
$ws.Range("M98").Value = 643.2
# Row 99
$ws.Range("H99").Value = 1418.2
$ws.Range("I99").Value = 1418.2
$ws.Range("K99").Value = 4254.6
$ws.Range("M99").Value = -2756.6
# Row 122
$ws.Range("H122").Value = 908.875
$ws.Range("I122").Value = 854.8
$ws.Range("K122").Value = 2564.4
$ws.Range("M122").Value = -114.3999999999996
# Row 125
$ws.Range("H125").Value = 7348
$ws.Range("I125").Value = 1319.8
$ws.Range("J125").Value = 12371.5
$ws.Range("K125").Value = 11878.2
$ws.Range("L125").Value = 111343.5
$ws.Range("M125").Value = -9418.199999999999
$ws.Range("N125").Value = -116263.5
# Row 131
$ws.Range("H131").Value = 16399.143
$ws.Range("I131").Value = 17114.834
$ws.Range("K131").Value = 51344.50199999999
$ws.Range("M131").Value = -46304.50199999999
# Row 137
$ws.Range("H137").Value = 1832.2222
$ws.Range("I137").Value = 1749.3334
$ws.Range("K137").Value = 5248.0002
$ws.Range("M137").Value = -2698.0002
# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 16500.2
$ws.Range("I63").Value = 19751
$ws.Range("K63").Value = 19751
$ws.Range("M63").Value = -19065
# Row 66
$ws.Range("H66").Value = 16500.2
$ws.Range("I66").Value = 19751
$ws.Range("K66").Value = 98755
$ws.Range("M66").Value = -95323

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 398.33334
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 97.5
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 97.5
$ws.Range("M5").Value = -887
$ws.Range("N5").Value = -323.5
# Row 86
$ws.Range("H86").Value = 1216
$ws.Range("I86").Value = 1185.1428
$ws.Range("J86").Value = 1259.2
$ws.Range("K86").Value = 1185.1428
$ws.Range("L86").Value = 1259.2
$ws.Range("M86").Value = -62.14280000000008
$ws.Range("N86").Value = -3505.2
# Row 89
$ws.Range("H89").Value = 1216
$ws.Range("I89").Value = 1185.1428
$ws.Range("J89").Value = 1259.2
$ws.Range("K89").Value = 5925.714
$ws.Range("L89").Value = 6296
$ws.Range("M89").Value = -309.7139999999999
$ws.Range("N89").Value = -17528
# Row 100
$ws.Range("H100").Value = 39828.43
$ws.Range("J100").Value = 39828.43
$ws.Range("L100").Value = 39828.43
$ws.Range("N100").Value = -41992.43
# Row 134
$ws.Range("H134").Value = 1581.4615
$ws.Range("I134").Value = 1581.4615
$ws.Range("K134").Value = 4744.3845
$ws.Range("M134").Value = -2209.3845

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 4071.5557
$ws.Range("J16").Value = 4177.7144
$ws.Range("L16").Value = 4177.7144
$ws.Range("N16").Value = -4751.7144
# Row 31
$ws.Range("H31").Value = 2437.5454
$ws.Range("I31").Value = 2437.5454
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2437.5454
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2142.5454
$ws.Range("N31").ClearContents()
# Row 34
$ws.Range("H34").Value = 2437.5454
$ws.Range("I34").Value = 2437.5454
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2437.5454
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2235.5454
$ws.Range("N34").ClearContents()
# Row 62
$ws.Range("H62").Value = 3333
$ws.Range("I62").Value = 3333
$ws.Range("K62").Value = 3333
$ws.Range("M62").Value = -2709
# Row 65
$ws.Range("H65").Value = 3333
$ws.Range("I65").Value = 3333
$ws.Range("K65").Value = 16665
$ws.Range("M65").Value = -13545
# Row 99
$ws.Range("H99").Value = 1964.5
$ws.Range("I99").Value = 1904.9333
$ws.Range("J99").Value = 2262.3333
$ws.Range("K99").Value = 1904.9333
$ws.Range("L99").Value = 2262.3333
$ws.Range("M99").Value = -406.9332999999999
$ws.Range("N99").Value = -5258.3333
# Row 113
$ws.Range("H113").Value = 4071.5557
$ws.Range("J113").Value = 4177.7144
$ws.Range("L113").Value = 4177.7144
$ws.Range("N113").Value = -8517.714400000001
# Row 126
$ws.Range("H126").Value = 1964.5
$ws.Range("I126").Value = 1904.9333
$ws.Range("J126").Value = 2262.3333
$ws.Range("K126").Value = 5714.7999
$ws.Range("L126").Value = 6786.999899999999
$ws.Range("M126").Value = -3244.7999
$ws.Range("N126").Value = -11726.9999

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 117
$ws.Range("H117").Value = 998
$ws.Range("I117").Value = 998
$ws.Range("K117").Value = 2994
$ws.Range("M117").Value = 448
# Row 121
$ws.Range("H121").Value = 10275.286
$ws.Range("I121").Value = 36808.332
$ws.Range("J121").Value = 3039
$ws.Range("K121").Value = 110424.996
$ws.Range("L121").Value = 9117
$ws.Range("M121").Value = -109114.996
$ws.Range("N121").Value = -11737
# Row 129
$ws.Range("H129").Value = 668256.1
$ws.Range("I129").Value = 1372.2222
$ws.Range("J129").Value = 1668582
$ws.Range("K129").Value = 4116.6666
$ws.Range("L129").Value = 5005746
$ws.Range("M129").Value = 883.3334000000004
$ws.Range("N129").Value = -5015746

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 1569.3
$ws.Range("I68").Value = 1115.8334
$ws.Range("K68").Value = 1115.8334
$ws.Range("M68").Value = -366.8334
# Row 71
$ws.Range("H71").Value = 1569.3
$ws.Range("I71").Value = 1115.8334
$ws.Range("K71").Value = 5579.166999999999
$ws.Range("M71").Value = -1835.166999999999
# Row 82
$ws.Range("H82").Value = 631.3333
$ws.Range("I82").Value = 432.57144
$ws.Range("J82").Value = 805.25
$ws.Range("K82").Value = 432.57144
$ws.Range("L82").Value = 805.25
$ws.Range("M82").Value = -71.57144
$ws.Range("N82").Value = -1527.25
# Row 85
$ws.Range("H85").Value = 631.3333
$ws.Range("I85").Value = 432.57144
$ws.Range("J85").Value = 805.25
$ws.Range("K85").Value = 432.57144
$ws.Range("L85").Value = 805.25
$ws.Range("M85").Value = 815.4285600000001
$ws.Range("N85").Value = -3301.25
# Row 93
$ws.Range("H93").Value = 1221.2307
$ws.Range("I93").Value = 1042.8889
$ws.Range("K93").Value = 1042.8889
$ws.Range("M93").Value = 205.1111000000001
# Row 103
$ws.Range("H103").Value = 18276.5
$ws.Range("J103").Value = 18276.5
$ws.Range("L103").Value = 18276.5
$ws.Range("N103").Value = -20620.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 4659.5
$ws.Range("I126").Value = 4898
$ws.Range("K126").Value = 14694
$ws.Range("M126").Value = -12224
